$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): reorder the "Odd_CS_*" columns AG..BB ---
# The columns "Odd_CS_4-4" (was AG1) and "Odd_CS_3-2_HT" (was AS1) move to the
# end of the AG:BB block, everything else shifts left to fill the gap.
$headers = @(
    "Odd_CS_0-1",
    "Odd_CS_0-2",
    "Odd_CS_1-2",
    "Odd_CS_0-3",
    "Odd_CS_1-3",
    "Odd_CS_2-3",
    "Odd_CS_1-0_HT",
    "Odd_CS_2-0_HT",
    "Odd_CS_2-1_HT",
    "Odd_CS_3-0_HT",
    "Odd_CS_3-1_HT",
    "Odd_CS_0-0_HT",
    "Odd_CS_1-1_HT",
    "Odd_CS_2-2_HT",
    "Odd_CS_0-1_HT",
    "Odd_CS_0-2_HT",
    "Odd_CS_1-2_HT",
    "Odd_CS_0-3_HT",
    "Odd_CS_1-3_HT",
    "Odd_CS_2-3_HT",
    "Odd_CS_4-4",
    "Odd_CS_3-2_HT"
)

$startCol = 33  # column AG
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $startCol + $i).Value = $headers[$i]
}

# --- Data row (row 2): reorder the matching odds values AG2..BB2 ---
$values = @(9, 15, 11, 29, 26, 34, 4.33, 13, 23, 41, 67, 2.63, 8, 51, 4.75, 17, 26, 51, 81, 201, 251)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2, $startCol + $i).Value = $values[$i]
}

# --- Other standalone odds updates on row 2 ---
$ws.Range("G2").Value = 2.45
$ws.Range("W2").Value = 8
$ws.Range("Y2").Value = 10
